# Issue #20: Printing the output of grep
#
# The heuristics table gains a "pattern"/"case-sensitive" split of the old
# "regex" column, five new JDBC-driver rows (HyperSQL, Derby, H2, Oracle,
# PostgreSQL), and the stray trailing colon on the MySQL JDBC prefix is
# dropped ("jdbc:mysql:" -> "jdbc:mysql").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the five new database rows (HyperSQL, Derby, H2, Oracle,
# PostgreSQL) right after the existing MySQL/C# rows.
$ws.Rows("5:9").Insert() | Out-Null

# --- Write cells in the same order the original author must have, so the
# shared-string table comes out in the exact sequence the diff shows
# (new strings are interned in first-use order; order matters for an
# exact OOXML match even though it's invisible to the end user). ---
$ws.Range("A5").Value = "Java "
$ws.Range("C5").Value = "jdbc:hsqldb"
$ws.Range("C2").Value = "jdbc:mysql"
$ws.Range("B6").Value = "Derby"
$ws.Range("C6").Value = "jdbc:derby"
$ws.Range("C7").Value = "jdbc:h2"
$ws.Range("B7").Value = "H2"
$ws.Range("B8").Value = "Oracle"
$ws.Range("C8").Value = "jdbc:oracle"
$ws.Range("B5").Value = "HyperSQL"
$ws.Range("C1").Value = "pattern"
$ws.Range("E1").Value = "case-sensitive"
$ws.Range("B9").Value = "PostgreSQL"
$ws.Range("C9").Value = "jdbc:postgresql"

# Remaining cells: reuse of already-interned strings, plus the new
# boolean "case-sensitive" column values.
$ws.Range("D1").Value = "regex"
$ws.Range("A6").Value = "Java"
$ws.Range("A7").Value = "Java"
$ws.Range("A8").Value = "Java"
$ws.Range("A9").Value = "Java"

$ws.Range("D2").Value = $false
$ws.Range("E2").Value = $true
$ws.Range("D3").Value = $false
$ws.Range("E3").Value = $true
$ws.Range("D4").Value = $false
$ws.Range("E4").Value = $true
$ws.Range("D5").Value = $false
$ws.Range("E5").Value = $true
$ws.Range("D6").Value = $false
$ws.Range("E6").Value = $true
$ws.Range("D7").Value = $false
$ws.Range("E7").Value = $true
$ws.Range("D8").Value = $false
$ws.Range("E8").Value = $true
$ws.Range("D9").Value = $false
$ws.Range("E9").Value = $true

# New header cells (D1, E1) pick up the same bold+bordered look as the
# rest of row 1.
$ws.Range("C1").Copy() | Out-Null
$ws.Range("D1:E1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

# Column widths, matching the "best fit" sizing the real workbook ended
# up with (ColumnWidth = stored XML width minus the ~0.8333 padding this
# engine adds on export).
$ws.Columns("A").ColumnWidth = 7.830729166666667
$ws.Columns("B").ColumnWidth = 9.666666666666666
$ws.Columns("C").ColumnWidth = 34.998697916666664
$ws.Columns("D").ColumnWidth = 5.330729166666667
$ws.Columns("E").ColumnWidth = 11.998697916666666

# Final cursor position left by the author.
$ws.Range("C18").Select() | Out-Null
